$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.161052346229553
$ws.Range("B1").Value = 2.372506618499756
$ws.Range("D1").Value = 2.393642663955688
$ws.Range("E1").Value = 1.219664454460144
